# update machinery info script
#
# 1. Remove the obsolete "Scripting" sheet entirely.
# 2. Lower-case the tab names for the main navigation sheets (this also
#    fixes up formulas on ModelOLD that pointed at "Main!..." -> "main!...").
# 3. Refresh the "use for scripting" flag row at the top of the Model sheet
#    so every quarter column is marked with an "x" (previously only I1 was),
#    and relabel the row header.

$wb = $excel.ActiveWorkbook

# --- 1. Delete the obsolete "Scripting" sheet -------------------------------
$wb.Worksheets.Item("Scripting").Delete() | Out-Null

# --- 2. Rename sheets --------------------------------------------------------
$wb.Worksheets.Item("Main").Name  = "main"
$wb.Worksheets.Item("Model").Name = "model"
$wb.Worksheets.Item("Parts").Name = "parts"
$wb.Worksheets.Item("Notes").Name = "notes"

# --- 3. Update the scripting-flag row on the Model sheet --------------------
$wsModel = $wb.Worksheets.Item("model")
$wsModel.Range("A1").Value    = "use for scripting"
$wsModel.Range("C1:I1").Value = "x"
$wsModel.Range("K1").Value    = "x"

# --- 4. Restore per-sheet view/selection state -------------------------------
$wsMain = $wb.Worksheets.Item("main")
$wsMain.Activate() | Out-Null
$wsMain.Range("B35").Select() | Out-Null

$wsParts = $wb.Worksheets.Item("parts")
$wsParts.Activate() | Out-Null
$wsParts.Range("C1").Select() | Out-Null

$wsNotes = $wb.Worksheets.Item("notes")
$wsNotes.Activate() | Out-Null
$wsNotes.Range("F36").Select() | Out-Null

$wsMfg = $wb.Worksheets.Item("Manufacturing")
$wsMfg.Activate() | Out-Null
$wsMfg.Range("C19").Select() | Out-Null
$excel.ActiveWindow.Zoom = 205

$wsModelOld = $wb.Worksheets.Item("ModelOLD")
$wsModelOld.Activate() | Out-Null
$wsModelOld.Range("C3").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 79
$wsModelOld.Range("C1").Select() | Out-Null
$wsModelOld.Range("A3").Select() | Out-Null
$wsModelOld.Range("F84").Select() | Out-Null

$wsIP = $wb.Worksheets.Item("IP")
$wsIP.Activate() | Out-Null
$wsIP.Range("B5").Select() | Out-Null
$excel.ActiveWindow.Zoom = 205

# "model" ends up the active/selected tab, matching the saved workbook state
$wsModel.Activate() | Out-Null
$wsModel.Range("B1").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 90
$wsModel.Range("C1:I1").Select() | Out-Null
